$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A / B (Hoje / Operadora) for the 3 new rows ---
$ws.Range("A3").Value = "13-04-2023"
$ws.Range("B3").Value = "417823 - PREMIUM SAÚDE S.A"
$ws.Range("A4").Value = "13-04-2023"
$ws.Range("B4").Value = "417823 - PREMIUM SAÚDE S.A"
$ws.Range("A5").Value = "13-04-2023"
$ws.Range("B5").Value = "417823 - PREMIUM SAÚDE S.A"

# --- Column C (Data da Notificação date-time) ---
$ws.Range("C2").Value = "03/04/2023  13:44:27"
$ws.Range("C3").Value = "03/04/2023  14:47:05"
$ws.Range("C4").Value = "03/04/2023  15:25:17"
$ws.Range("C5").Value = "03/04/2023  19:51:17"

# --- Column D / E (Demanda / Protocolo numbers) ---
$ws.Range("D2").Value = 12153205
$ws.Range("E2").Value = 8588077
$ws.Range("D3").Value = 12153463
$ws.Range("E3").Value = 8588431
$ws.Range("D4").Value = 12153630
$ws.Range("E4").Value = 8588645
$ws.Range("D5").Value = 12154294
$ws.Range("E5").Value = 8589407

# --- Column F (Beneficiário) ---
$ws.Range("F2").Value = "BÁRBARA DE OLIVEIRA PATRÍCIO"
$ws.Range("F3").Value = "KEPA FREDRICK STOCKNER"
$ws.Range("F4").Value = "DAVIDSON LUIZ PEREIRA LOPES"
$ws.Range("F5").Value = "KENIA CLAUDIA FARIA CAMPOS"

# --- Column G (Prazo) ---
$ws.Range("G2").Value = "4 dias úteis"
$ws.Range("G3").Value = "4 dias úteis"
$ws.Range("G4").Value = "4 dias úteis"
$ws.Range("G5").Value = "4 dias úteis"

# --- Column H (Respondido) ---
$ws.Range("H2").Value = "NO"
$ws.Range("H3").Value = "NO"
$ws.Range("H4").Value = "NO"
$ws.Range("H5").Value = "NO"

# --- Column I (Natureza) ---
$ws.Range("I2").Value = "Assistencial"
$ws.Range("I3").Value = "Assistencial"
$ws.Range("I4").Value = "Assistencial"
$ws.Range("I5").Value = "Assistencial"

# --- Column J (Opções) ---
$ws.Range("J2").Value = "Responder  Detalhes"
$ws.Range("J3").Value = "Responder  Detalhes"
$ws.Range("J4").Value = "Responder  Detalhes"
$ws.Range("J5").Value = "Responder  Detalhes"
